$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the "receiving_uom_code" column (was C, now E) to hold
# the new "brand name" / "manufacture by" fields.
$ws.Columns("C:D").Insert()

# Header row (row 1)
$ws.Range("C1").Value = "brand name"
$ws.Range("D1").Value = "manufacture by"

# Data row (row 2)
$ws.Range("C2").Value = "Amul"
$ws.Range("D2").Value = "Amul"

# Column widths (best-fit-like sizing for the new columns; existing columns B, E, F, G
# keep the widths they already had, shifted along automatically by the column insert).
$ws.Columns("A").ColumnWidth = 5.16
$ws.Columns("C").ColumnWidth = 11.66
$ws.Columns("D").ColumnWidth = 16.5

# The active selection moves from C7 to D7.
$ws.Range("D7").Select() | Out-Null
